$d = $word.ActiveDocument

# 1. Merge the "Save password encrypted " + bookmark + "using bcrypt..." runs
#    into a single run with the full combined text, and drop the now-redundant
#    _GoBack bookmark that used to sit in the middle of that sentence.
$d.Content.Find.Execute(
    "Save password encrypted using bcrypt lib from spring security with 12 round.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Save password encrypted using bcrypt lib from spring security with 12 round.",
    2) | Out-Null

# 2. Append five new paragraphs at the end of the document describing the
#    "Fifth day" tasks, matching the original formatting conventions used
#    throughout the document (bullet list numId=1 for list items, numId=0 /
#    leftChars=0 for the blank spacer paragraphs, and a plain paragraph with
#    no numbering for the day heading).
function New-FlatOpcPackage($innerBodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerBodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

function Insert-ParagraphsAtEnd($innerBodyXml) {
    $pos = $d.Content.End
    $r = $d.Range($pos, $pos)
    $r.InsertXML((New-FlatOpcPackage $innerBodyXml)) | Out-Null
}

$rPr = '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>'
$bulletPPr = '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="420" w:leftChars="0" w:hanging="420" w:firstLineChars="0"/>' + $rPr + '</w:pPr>'
$blankPPr = '<w:pPr><w:numPr><w:numId w:val="0"/></w:numPr><w:ind w:leftChars="0"/>' + $rPr + '</w:pPr>'
$plainPPr = '<w:pPr>' + $rPr + '</w:pPr>'

$pBlank1 = '<w:p>' + $blankPPr + '</w:p>'
$pHeading = '<w:p>' + $plainPPr + '<w:r>' + $rPr + '<w:t>Fifth day(monday):</w:t></w:r></w:p>'
$pSetup = '<w:p>' + $bulletPPr + '<w:r>' + $rPr + '<w:t>Setup security file.</w:t></w:r></w:p>'
$pJwt = '<w:p>' + $bulletPPr + '<w:r>' + $rPr + '<w:t>Implement jwt validation.</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$pBlank2 = '<w:p>' + $blankPPr + '</w:p>'

Insert-ParagraphsAtEnd ($pBlank1 + $pHeading + $pSetup + $pJwt + $pBlank2)
